# Show the list role/function next to each member in the overview export:
# the "Members" data cell used to just dump the member list; switch it to
# render each member together with their list role/function.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = '${twig:record.getListMembersWithFunctions()}'

# Cosmetic: leave the selection where the author's Excel session ended up.
$ws.Range("D9").Select() | Out-Null
